# EvaluacionPorCompetenciasTemplate.xlsx edit script
# Summary of the change:
#  - The evaluation table (header row 31, data rows 32-41) is moved two
#    columns to the left (from E:J to C:H) and a new column
#    "Calificación VE" is inserted into it (yielding C:I).
#  - Two headers are renamed: "Column1" -> "Varible equivalente (VE)"
#    and "Variable equivalente" -> "Método seleccionado".
#  - A "Promedio ponderado" total value is filled in (G44).
#  - Both charts' series formulas are repointed at the new ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------
# 1. Move the table two columns to the left: E31:J41 -> C31:H41
#    (staged through a far-away scratch range so the cut source and
#    destination never overlap, which would otherwise drop values)
# ---------------------------------------------------------------
$ws.Range("E31:J41").Cut($ws.Range("AA31"))
$ws.Range("AA31:AF41").Cut($ws.Range("C31"))
$tbl.Resize($ws.Range("C31:H41"))

# ---------------------------------------------------------------
# 2. Insert a new column for "Calificación VE" between the
#    "Método seleccionado" column and "Valor Relativo" (only within
#    the table's rows, so nothing else on the sheet shifts)
# ---------------------------------------------------------------
$ws.Range("F31:F41").Insert(-4161)
$tbl.Resize($ws.Range("C31:I41"))

# Clean up the cells vacated by the shift/insert dance
$ws.Range("J31:K41").Clear()
$ws.Range("AA31:AG41").Clear()

# ---------------------------------------------------------------
# 3. Header text
# ---------------------------------------------------------------
$ws.Range("C31").Value = "Méritos y Habilidades"
$ws.Range("D31").Value = "Varible equivalente (VE)"
$ws.Range("E31").Value = "Método seleccionado"
$ws.Range("F31").Value = "Calificación VE"
$ws.Range("G31").Value = "Valor Relativo"
$ws.Range("H31").Value = "Candidato Ideal"
$ws.Range("I31").Value = "Calificación"

for ($i = 1; $i -le $tbl.ListColumns.Count; $i++) {
  $col = $tbl.ListColumns.Item($i)
}

# ---------------------------------------------------------------
# 4. Promedio ponderado total
# ---------------------------------------------------------------
$ws.Range("G44").Value = 2.5875498499999998
$ws.Rows("44:44").RowHeight = 21

# ---------------------------------------------------------------
# 5. Merge E6:F6 grows to E6:G6
# ---------------------------------------------------------------
$ws.Range("E6:F6").UnMerge()
$ws.Range("E6:G6").Merge()

# ---------------------------------------------------------------
# 6. Repoint both charts at the new ranges
# ---------------------------------------------------------------
for ($c = 1; $c -le 2; $c++) {
  $co = $ws.ChartObjects($c)
  $ch = $co.Chart
  $s1 = $ch.SeriesCollection(1)
  $s1.Formula = "=SERIES(Hoja1!`$H`$31,Hoja1!`$C`$32:`$C`$41,Hoja1!`$H`$32:`$H`$41,1)"
  $s2 = $ch.SeriesCollection(2)
  $s2.Formula = "=SERIES(Hoja1!`$I`$31,Hoja1!`$C`$32:`$C`$41,Hoja1!`$I`$32:`$I`$41,2)"
}
